$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.315045833587646
$ws.Range("B1").Value = 1.755271911621094
$ws.Range("C1").Value = 4.193210601806641
$ws.Range("D1").Value = 3.043334484100342
$ws.Range("E1").Value = 1.123455405235291
